$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "50.709.77"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "2.901.31"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "369.68"
$ws.Range("E5").Value = "  +4.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.84"
$ws.Range("E6").Value = "  -4.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("E7").Value = "  -3.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.03"
$ws.Range("E8").Value = "  +3.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -4.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.85"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0834"
$ws.Range("E12").Value = "  -3.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.42"
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("D14").Value = "3.352.98"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.37"
$ws.Range("E15").Value = "  -4.60%  "
$ws.Range("D16").Value = "2.891.15"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.937"
$ws.Range("E17").Value = "  -2.97%  "
$ws.Range("D18").Value = "50.762.25"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.26"
$ws.Range("E19").Value = "  -5.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.21"
$ws.Range("E20").Value = "  -3.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.80"
$ws.Range("E21").Value = "  -5.26%  "
$ws.Range("D22").Value = "0.0₃0943"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.19"
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "259.60"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.69"
$ws.Range("E25").Value = "  -2.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.171"
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.70"
$ws.Range("E28").Value = "  -4.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.99"
$ws.Range("E29").Value = "  -7.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.102"
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.87"
$ws.Range("E31").Value = "  -3.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.08"
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.13"
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.63"
$ws.Range("E34").Value = "  -5.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.03"
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0418"
$ws.Range("E37").Value = "  -2.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.04"
$ws.Range("E38").Value = "  -3.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.63"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.02"
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.84"
$ws.Range("E41").Value = "  -6.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.113"
$ws.Range("E42").Value = "  -3.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.18"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.06"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.08"
$ws.Range("E45").Value = "  -3.69%  "
$ws.Range("D46").Value = "2.014.41"
$ws.Range("E46").Value = "  -4.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.17"
$ws.Range("E48").Value = "  -5.62%  "
$ws.Range("D49").Value = "3.180.64"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.235"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0311"
$ws.Range("E51").Value = "  -8.66%  "
